$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4-6 Pour Over Recipe")

# --- Brew ratio input changed from 15 to 17 ---
$ws.Range("C5").Value = 17

# --- Header row 13: insert a new "Total (g)" column before "Leaves (g)" ---
# Shift the existing "Leaves (g)" header (with its formatting) from E13 to F13,
# then put the new "Total (g)" label in E13 (same formatting it already had).
$ws.Range("E13").Copy($ws.Range("F13"))
$ws.Range("E13").Value = "Total (g)"

# --- Rows 14-19: insert a new column E that is a running total of water poured,
#     pushing the previous "water remaining" column from E to F ---
$ws.Range("E14").Copy($ws.Range("F14"))
$ws.Range("E15").Copy($ws.Range("F15"))
$ws.Range("E16").Copy($ws.Range("F16"))
$ws.Range("E17").Copy($ws.Range("F17"))
$ws.Range("E18").Copy($ws.Range("F18"))
$ws.Range("E19").Copy($ws.Range("F19"))

# New running-total formulas in column E
$ws.Range("E14").Formula = "=D14"
$ws.Range("E15").Formula = "=E14+D15"
$ws.Range("E16:E19").Formula = "=E15+D16"

# Fix up the "water remaining" formulas that moved to column F so they reference
# the new column-F cells instead of the old column-E cells
$ws.Range("F14").Formula = "=water-bloom"
$ws.Range("F15").Formula = "=F14-first_pour"
$ws.Range("F16").Formula = "=F15-second_pour"
$ws.Range("F17").Formula = "=F16-third_pour"
$ws.Range("F18").Formula = "=F17-fourth_pour"
$ws.Range("F19").Formula = "=F18-fifth_pour"

# --- Shift the blank styled filler cells in column G to column H for rows 14-19,
#     and G13 to K13 (row 13 already has H/I/J filled) ---
$ws.Range("G14").Copy($ws.Range("H14"))
$ws.Range("G14").Clear()
$ws.Range("G15").Copy($ws.Range("H15"))
$ws.Range("G15").Clear()
$ws.Range("G16").Copy($ws.Range("H16"))
$ws.Range("G16").Clear()
$ws.Range("G17").Copy($ws.Range("H17"))
$ws.Range("G17").Clear()
$ws.Range("G18").Copy($ws.Range("H18"))
$ws.Range("G18").Clear()
$ws.Range("G19").Copy($ws.Range("H19"))
$ws.Range("G19").Clear()

$ws.Range("G13").Copy($ws.Range("K13"))
$ws.Range("G13").Clear()

# Restore the selection to match what was last saved
$ws.Range("B1").Select()
